$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: update title (D9) and link (E9)
$ws.Range("D9").Value = "왜 공대와 경영학과는 저렇게 생각할까?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/depth-level-angle-knowledge/#utm_source=rss&utm_medium=rss&utm_campaign=depth-level-angle-knowledge"

# Row 12: update title (D12) and link (E12)
$ws.Range("D12").Value = "“개발자를 위한 AI & ML” 도서 번역을 마쳤습니다."
$ws.Range("E12").Value = "https://tensorflow.blog/2022/07/13/%ea%b0%9c%eb%b0%9c%ec%9e%90%eb%a5%bc-%ec%9c%84%ed%95%9c-ai-ml-%eb%8f%84%ec%84%9c-%eb%b2%88%ec%97%ad%ec%9d%84-%eb%a7%88%ec%b3%a4%ec%8a%b5%eb%8b%88%eb%8b%a4/"

# Row 51: update title (D51) and link (E51)
$ws.Range("D51").Value = "[MySQL] 데이터베이스 생성하는 방법"
$ws.Range("E51").Value = "https://bskyvision.com/entry/MySQL-%EB%8D%B0%EC%9D%B4%ED%84%B0%EB%B2%A0%EC%9D%B4%EC%8A%A4-%EC%83%9D%EC%84%B1%ED%95%98%EB%8A%94-%EB%B0%A9%EB%B2%95"
